$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header cell H1 "Save" with the same style as the other header cells (copy from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new data values in column H for rows 2 and 3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
